# B1--and-B2-PowerPoint.pptx — Tue, Apr 21, 2020  3:06:28 PM
#
# 1) The table on slide 5 gets a new table style (its StyleId/GUID changes).
# 2) The deck's theme colour palette is switched from the custom
#    "Integral" (Red Violet) scheme to the stock Office palette
#    (the same RGB values the notes-master's "Office Theme" theme part
#    already used).

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 5 -------------------------------------
$newTableStyleId = "{DE8CE84C-CED6-4BD7-8179-E3F952A70CA8}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# --- 2. Swap the theme colour scheme over to the stock Office palette ----
function Convert-HexToRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Slot order matches MsoThemeColorSchemeIndex / ppColorSchemeIndex:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$officeColors = [ordered]@{
    1  = "000000"   # dk1
    2  = "FFFFFF"   # lt1
    3  = "44546A"   # dk2
    4  = "E7E6E6"   # lt2
    5  = "5B9BD5"   # accent1
    6  = "ED7D31"   # accent2
    7  = "A5A5A5"   # accent3
    8  = "FFC000"   # accent4
    9  = "4472C4"   # accent5
    10 = "70AD47"   # accent6
    11 = "0563C1"   # hlink
    12 = "954F72"   # folHlink
}

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
foreach ($slot in $officeColors.Keys) {
    $themeColors.Item($slot).RGB = Convert-HexToRgb $officeColors[$slot]
}
